# Update Pooh Points site
# Applies the latest box-score snapshot (game clock 3:51 - 1st Half) to the
# "Players" sheet and the resulting starter totals to the "OwnerTotals" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Players")
$ws2 = $wb.Worksheets.Item("OwnerTotals")

# --- Players sheet (row 2) ---
$ws1.Cells.Item(2, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(2, 8).Value = 5
$ws1.Cells.Item(2, 15).Value = 12

# --- Players sheet (row 3) ---
$ws1.Cells.Item(3, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(3, 8).Value = 18
$ws1.Cells.Item(3, 9).Value = 14
$ws1.Cells.Item(3, 11).Value = 4
$ws1.Cells.Item(3, 15).Value = 14

# --- Players sheet (row 4) - now Denzel Aberdeen (UK) ---
$ws1.Cells.Item(4, 4).Value = 'Denzel Aberdeen'
$ws1.Cells.Item(4, 5).Value = 'UK'
$ws1.Cells.Item(4, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(4, 8).Value = 4
$ws1.Cells.Item(4, 9).Value = 4
$ws1.Cells.Item(4, 10).Value = 0
$ws1.Cells.Item(4, 15).Value = 7

# --- Players sheet (row 5) - now Felix Okpara (TENN) ---
$ws1.Cells.Item(5, 4).Value = 'Felix Okpara'
$ws1.Cells.Item(5, 5).Value = 'TENN'
$ws1.Cells.Item(5, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(5, 8).Value = 2
$ws1.Cells.Item(5, 9).Value = 3
$ws1.Cells.Item(5, 10).Value = 1
$ws1.Cells.Item(5, 15).Value = 10

# --- Players sheet (row 6) ---
$ws1.Cells.Item(6, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(6, 10).Value = 3
$ws1.Cells.Item(6, 15).Value = 8

# --- Players sheet (row 7) ---
$ws1.Cells.Item(7, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(7, 8).Value = -1
$ws1.Cells.Item(7, 14).Value = 3
$ws1.Cells.Item(7, 15).Value = 12

# --- Players sheet (row 8) ---
$ws1.Cells.Item(8, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(8, 15).Value = 10

# --- Players sheet (row 9) ---
$ws1.Cells.Item(9, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(9, 8).Value = 1
$ws1.Cells.Item(9, 10).Value = 2
$ws1.Cells.Item(9, 15).Value = 11

# --- Players sheet (row 10) ---
$ws1.Cells.Item(10, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(10, 12).Value = 1
$ws1.Cells.Item(10, 14).Value = 1
$ws1.Cells.Item(10, 15).Value = 6

# --- Players sheet (row 11) ---
$ws1.Cells.Item(11, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(11, 8).Value = 3
$ws1.Cells.Item(11, 12).Value = 1
$ws1.Cells.Item(11, 15).Value = 6

# --- Players sheet (row 12) - now DeWayne Brown II (TENN) ---
$ws1.Cells.Item(12, 4).Value = 'DeWayne Brown II'
$ws1.Cells.Item(12, 5).Value = 'TENN'
$ws1.Cells.Item(12, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(12, 9).Value = 5
$ws1.Cells.Item(12, 10).Value = 1
$ws1.Cells.Item(12, 12).Value = 1
$ws1.Cells.Item(12, 13).Value = 1
$ws1.Cells.Item(12, 14).Value = 0
$ws1.Cells.Item(12, 15).Value = 7

# --- Players sheet (row 13) - now Jasper Johnson (UK) ---
$ws1.Cells.Item(13, 4).Value = 'Jasper Johnson'
$ws1.Cells.Item(13, 5).Value = 'UK'
$ws1.Cells.Item(13, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(13, 8).Value = 8
$ws1.Cells.Item(13, 9).Value = 8
$ws1.Cells.Item(13, 10).Value = 0
$ws1.Cells.Item(13, 12).Value = 0
$ws1.Cells.Item(13, 13).Value = 0
$ws1.Cells.Item(13, 14).Value = 1
$ws1.Cells.Item(13, 15).Value = 9

# --- Players sheet (row 14) ---
$ws1.Cells.Item(14, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(14, 10).Value = 3
$ws1.Cells.Item(14, 15).Value = 8

# --- Players sheet (rows 15-17) - clock only ---
$ws1.Cells.Item(15, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(16, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(17, 7).Value = '3:51 - 1st Half'

# --- Players sheet (row 18) - now Troy Henderson ---
$ws1.Cells.Item(18, 4).Value = 'Troy Henderson'
$ws1.Cells.Item(18, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(18, 11).Value = 1
$ws1.Cells.Item(18, 15).Value = 3

# --- Players sheet (row 19) - now Amaree Abram ---
$ws1.Cells.Item(19, 4).Value = 'Amaree Abram'
$ws1.Cells.Item(19, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(19, 8).Value = -1
$ws1.Cells.Item(19, 11).Value = 0
$ws1.Cells.Item(19, 15).Value = 2

# --- Players sheet (rows 20-21) - clock only ---
$ws1.Cells.Item(20, 7).Value = '3:51 - 1st Half'
$ws1.Cells.Item(21, 7).Value = '3:51 - 1st Half'

# --- OwnerTotals sheet ---
$ws2.Cells.Item(2, 2).Value = 18
$ws2.Cells.Item(3, 2).Value = 5
$ws2.Cells.Item(4, 1).Value = 'Mark'
$ws2.Cells.Item(6, 1).Value = 'Ron'
$ws2.Cells.Item(6, 3).Value = 0
$ws2.Cells.Item(7, 1).Value = 'Hal'
$ws2.Cells.Item(7, 2).Value = -1
$ws2.Cells.Item(7, 3).Value = 1
